# Update the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) timestamps on the zh-cn and de-de
# report sheets to reflect the newly generated handback report times.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 23:04:03"
$wsZhCn.Range("H2").Value = "2016-03-21 23:04:24"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 23:04:08"
$wsDeDe.Range("H2").Value = "2016-03-21 23:04:30"
